# "wrapping up test file audit"
#
# The "optimization_parameters" sheet had a stray leftover row (A16="Sheet",
# B16=3, C16=4) that doesn't belong with the real data. Remove it so the
# "simulation_timepoints" row shifts up from row 17 to row 16.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows.Item(16).Delete() | Out-Null

# Keep the (now shifted) simulation_timepoints row selected, same as the
# row that used to be selected before the delete.
$ws.Rows.Item(16).Select() | Out-Null

# While auditing, a cell was also clicked on the network_weights sheet.
$wb.Worksheets.Item("network_weights").Range("E10").Select() | Out-Null

# Finish on the threshold_b sheet/tab.
$wb.Worksheets.Item("threshold_b").Activate() | Out-Null
